$wb = $excel.ActiveWorkbook

# 1. Rename the "Include from EntityNamePartQu" sheet to "Include #0"
$wsInclude = $wb.Worksheets.Item("Include from EntityNamePartQu")
$wsInclude.Name = "Include #0"

# 2. Update the Metadata sheet
$ws = $wb.Worksheets.Item("Metadata")

# 2a. Update Version value
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# 2b. Update Date value
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# 2c. Insert a new "Jurisdiction" row after row 10 (Contact), shifting remaining rows down.
$ws.Range("A11:B11").Copy()
$ws.Rows.Item(11).EntireRow.Insert()
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Fix up formatting on the newly inserted row so it matches the body style used
# by the rest of the table (the raw row-insert above leaves a slightly different style).
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
